# update code tinh luong
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Đơn phụ phẫu 1": insert 2 new detail rows (2 new HD-LUXURY
# lines) above the "Tổng" row, then refresh the "Tổng" row's totals.
# -----------------------------------------------------------------
$wsPP1 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Push the existing "Tổng" row (row 12) down by two rows so we can
# insert the two new data rows above it.
$wsPP1.Rows.Item(12).Insert()
$wsPP1.Rows.Item(12).Insert()

# Column C holds the "Ngày thực hiện" dates as plain text (e.g.
# "07-28-2024"), matching every other row in this column. Format the
# two new cells as Text first so Excel doesn't silently reinterpret
# the typed value as a real date.
$wsPP1.Range("C12:C13").NumberFormat = "@"

# New data row 12
$wsPP1.Cells.Item(12, 1).Value2 = "HD-LUXURY"
$wsPP1.Cells.Item(12, 2).Value2 = 593
$wsPP1.Cells.Item(12, 3).Value2 = "07-28-2024"
$wsPP1.Cells.Item(12, 4).Value2 = "SÓC TRĂNG"
$wsPP1.Cells.Item(12, 5).Value2 = "trần thị ngọc bích "
$wsPP1.Cells.Item(12, 6).Value2 = "Cá nhân"
$wsPP1.Cells.Item(12, 7).Value2 = "Nâng mũi"
$wsPP1.Cells.Item(12, 8).Value2 = "Kha Như Huỳnh "
$wsPP1.Cells.Item(12, 9).Value2 = 100000

# New data row 13
$wsPP1.Cells.Item(13, 1).Value2 = "HD-LUXURY"
$wsPP1.Cells.Item(13, 2).Value2 = 594
$wsPP1.Cells.Item(13, 3).Value2 = "07-28-2024"
$wsPP1.Cells.Item(13, 4).Value2 = "SÓC TRĂNG"
$wsPP1.Cells.Item(13, 5).Value2 = "thạch thị sơ ri"
$wsPP1.Cells.Item(13, 6).Value2 = "Cá nhân"
$wsPP1.Cells.Item(13, 7).Value2 = "Nâng mũi"
$wsPP1.Cells.Item(13, 8).Value2 = "Kha Như Huỳnh "
$wsPP1.Cells.Item(13, 9).Value2 = 100000

# Refresh the "Tổng" row, now shifted to row 14
$wsPP1.Cells.Item(14, 1).Value2 = "Tổng"
$wsPP1.Cells.Item(14, 2).Value2 = 12
$wsPP1.Cells.Item(14, 9).Value2 = 900000

# -----------------------------------------------------------------
# Sheet "Lương": the per-branch summary gained a "Tổng công tại ..."
# and "Phụ cấp tại ..." row for each branch (CẦN THƠ / LONG XUYÊN /
# SÓC TRĂNG), plus a new "Lương công tác tại ..." row for CẦN THƠ and
# LONG XUYÊN. This re-numbers everything below row 1. Easiest to
# rebuild rows 2-40 (previously 2-34) from scratch with the final
# labels/values taken from the updated report.
# -----------------------------------------------------------------
$wsL = $wb.Worksheets.Item("Lương")

$rows = @(
    @("Tổng công tại CẦN THƠ", 0),
    @("Phụ cấp tại CẦN THƠ", 0),
    @("Lương công tác tại CẦN THƠ", 0),
    @("Lương cơ bản tại CẦN THƠ", $null),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Chiết khấu thu nợ tại CẦN THƠ", 0),
    @("Ứng lương tại CẦN THƠ", -0),
    @("Tổng công tại LONG XUYÊN", 0),
    @("Phụ cấp tại LONG XUYÊN", 0),
    @("Lương công tác tại LONG XUYÊN", 0),
    @("Lương cơ bản tại LONG XUYÊN", $null),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Chiết khấu thu nợ tại LONG XUYÊN", 0),
    @("Ứng lương tại LONG XUYÊN", -0),
    @("Tổng công tại SÓC TRĂNG", 23),
    @("Phụ cấp tại SÓC TRĂNG", 805000),
    @("Lương cơ bản tại SÓC TRĂNG", 4653392.857142857),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 240000),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 16000),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 428800),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 900000),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0),
    @("Chiết khấu thu nợ tại SÓC TRĂNG", 100000),
    @("Ứng lương tại SÓC TRĂNG", -2000000),
    @("Tổng lương tại CẦN THƠ", 0),
    @("Tổng lương tại LONG XUYÊN", 0),
    @("Tổng lương tại SÓC TRĂNG", 5143192.857142857),
    @("Tổng lương", 5143192.857142857)
)

# The table grows from 33 data rows (old rows 2-34) to 39 data rows
# (new rows 2-40), so make room for 6 more rows before rewriting the
# label/value pairs in order.
for ($i = 0; $i -lt 6; $i++) {
    $wsL.Rows.Item(2).Insert()
}

$r = 2
foreach ($pair in $rows) {
    $wsL.Cells.Item($r, 1).Value2 = $pair[0]
    if ($null -eq $pair[1]) {
        $wsL.Cells.Item($r, 2).Value2 = $null
    } else {
        $wsL.Cells.Item($r, 2).Value2 = $pair[1]
    }
    $r = $r + 1
}
